$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 used to be "cocas" / 7.5 / 5 -- the author overtyped it with
# "borrado" placeholders in all three columns instead of removing the row.
$ws.Range("A5").Value = "borrado"
$ws.Range("B5").Value = "borrado"
$ws.Range("C5").Value = "borrado"

# A new product row was appended at the bottom of the table.
$ws.Range("A8").Value = "roles de canela"
$ws.Range("B8").Value = 15
$ws.Range("C8").Value = 4

# The author finished by clicking into the next empty row below the table.
[void]$ws.Range("B9").Select()
